{"js": "// Replace each two-digit multiplication expression in the document's\n// table cells with its updated counterpart, per the commit diff.\n// Each old expression is unique within the document, so an exact\n// (case-sensitive, non-wildcard) search safely targets exactly one run.\nconst replacements = [\n  [\"19\u00d732=\", \"64\u00d784=\"],\n  [\"50\u00d744=\", \"63\u00d758=\"],\n  [\"62\u00d754=\", \"91\u00d723=\"],\n  [\"36\u00d725=\", \"14\u00d755=\"],\n  [\"20\u00d760=\", \"63\u00d738=\"],\n  [\"64\u00d748=\", \"18\u00d713=\"],\n  [\"94\u00d735=\", \"96\u00d782=\"],\n  [\"68\u00d783=\", \"43\u00d714=\"],\n  [\"70\u00d731=\", \"46\u00d731=\"],\n  [\"90\u00d729=\", \"93\u00d747=\"],\n  [\"20\u00d714=\", \"72\u00d798=\"],\n  [\"23\u00d751=\", \"97\u00d788=\"],\n  [\"71\u00d736=\", \"78\u00d768=\"],\n  [\"26\u00d794=\", \"76\u00d747=\"],\n  [\"22\u00d794=\", \"17\u00d763=\"],\n  [\"69\u00d767=\", \"62\u00d741=\"],\n  [\"84\u00d731=\", \"59\u00d725=\"],\n  [\"53\u00d766=\", \"93\u00d728=\"],\n  [\"60\u00d738=\", \"30\u00d766=\"],\n  [\"16\u00d785=\", \"52\u00d713=\"],\n  [\"85\u00d781=\", \"27\u00d764=\"],\n  [\"99\u00d786=\", \"64\u00d735=\"],\n  [\"36\u00d763=\", \"13\u00d716=\"],\n  [\"18\u00d788=\", \"53\u00d779=\"],\n  [\"11\u00d754=\", \"55\u00d783=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication expression in the document's\n# table cells with its updated counterpart, per the commit diff.\n# Each old expression occurs exactly once in the document, so a single\n# Find/Replace (wdReplaceOne) per pair safely retargets just that run.\n\n$pairs = @(\n    @(\"19\u00d732=\", \"64\u00d784=\"),\n    @(\"50\u00d744=\", \"63\u00d758=\"),\n    @(\"62\u00d754=\", \"91\u00d723=\"),\n    @(\"36\u00d725=\", \"14\u00d755=\"),\n    @(\"20\u00d760=\", \"63\u00d738=\"),\n    @(\"64\u00d748=\", \"18\u00d713=\"),\n    @(\"94\u00d735=\", \"96\u00d782=\"),\n    @(\"68\u00d783=\", \"43\u00d714=\"),\n    @(\"70\u00d731=\", \"46\u00d731=\"),\n    @(\"90\u00d729=\", \"93\u00d747=\"),\n    @(\"20\u00d714=\", \"72\u00d798=\"),\n    @(\"23\u00d751=\", \"97\u00d788=\"),\n    @(\"71\u00d736=\", \"78\u00d768=\"),\n    @(\"26\u00d794=\", \"76\u00d747=\"),\n    @(\"22\u00d794=\", \"17\u00d763=\"),\n    @(\"69\u00d767=\", \"62\u00d741=\"),\n    @(\"84\u00d731=\", \"59\u00d725=\"),\n    @(\"53\u00d766=\", \"93\u00d728=\"),\n    @(\"60\u00d738=\", \"30\u00d766=\"),\n    @(\"16\u00d785=\", \"52\u00d713=\"),\n    @(\"85\u00d781=\", \"27\u00d764=\"),\n    @(\"99\u00d786=\", \"64\u00d735=\"),\n    @(\"36\u00d763=\", \"13\u00d716=\"),\n    @(\"18\u00d788=\", \"53\u00d779=\"),\n    @(\"11\u00d754=\", \"55\u00d783=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $matched = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $matched) {\n        Write-Output \"WARNING: no match found for $oldText\"\n    }\n}\n\n$d.Saved = $false\n"}
